$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 520; this shifts the existing rows 520:645 down
# to 521:646 (and the sheet dimension grows to A1:R646), reproducing the
# effect of a new weekly price record being prepended before the old row 520.
$ws.Rows.Item(520).Insert()

# Populate the new row 520 with a new weekly record: identical to the
# (now shifted-down) row 521 data except for the reporting date (column D),
# which moves forward one week.
$ws.Range("A520").Value = 9
$ws.Range("B520").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C520").Value = "Metropolitana"
$ws.Range("D520").Value = 45135
$ws.Range("E520").Value = 13
$ws.Range("F520").Value = 100112012
$ws.Range("G520").Value = "Espinaca"
$ws.Range("H520").Value = "Sin especificar"
$ws.Range("I520").Value = "Primera"
$ws.Range("J520").Value = 160
$ws.Range("K520").Value = 6000
$ws.Range("L520").Value = 8000
$ws.Range("M520").Value = 7000
$ws.Range("N520").Value = "`$/cuna 10 kilos"
$ws.Range("O520").Value = "Provincia de Chacabuco"
$ws.Range("P520").Value = 700
$ws.Range("Q520").Value = 10
$ws.Range("R520").Value = "Hortaliza"
